$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The rows being appended (8-13) duplicate the existing data rows
# (2,7,4,5,3,6) in this specific order.
$sourceRows = @(2, 7, 4, 5, 3, 6)

$destRow = 8
foreach ($src in $sourceRows) {
    $ws.Range("A$src`:K$src").Copy()
    $ws.Range("A$destRow`:K$destRow").PasteSpecial()
    $destRow++
}
